# Commit: "add Actions class comment"
#
# The underlying data changes in this revision:
#  - LoginTest sheet (sheet2): credentials updated to bluebirdsr@pobox.sk /
#    Lampa2019!; mailto: hyperlinks removed (plain text now); sheet becomes
#    the active tab/selection.
#  - CreateAccountTest sheet (sheet3): previously held only an
#    "accountname"/"Raman" pair; it is rebuilt with the same
#    username/password/runmode layout as LoginTest (same credential values),
#    dropping the old accountname/Raman data entirely.
#  - Column widths / selections on both sheets adjusted to match.

$wb = $excel.ActiveWorkbook

# ---- LoginTest sheet -------------------------------------------------
$wsLogin = $wb.Worksheets.Item("LoginTest")

# Drop the mailto: hyperlinks (and their relationships) - values stay, the
# cells keep their existing "Hyperlink" cell style.
[void]$wsLogin.Hyperlinks.Delete()

$wsLogin.Range("A2").Value = "bluebirdsr@pobox.sk"
$wsLogin.Range("B2").Value = "Lampa2019!"

$wsLogin.Columns.Item(1).ColumnWidth = 22
$wsLogin.Columns.Item(2).ColumnWidth = 24.333333333333332

# ---- CreateAccountTest sheet ------------------------------------------
$wsCreate = $wb.Worksheets.Item("CreateAccountTest")

# Wipe the old single-column accountname/Raman data ...
$wsCreate.Cells.Clear()

# ... and replace it with the same username/password/runmode table as
# LoginTest.
$wsCreate.Range("A1").Value = "username"
$wsCreate.Range("B1").Value = "password"
$wsCreate.Range("C1").Value = "runmode"
$wsCreate.Range("A2").Value = "bluebirdsr@pobox.sk"
$wsCreate.Range("B2").Value = "Lampa2019!"
$wsCreate.Range("C2").Value = "Y"

$wsCreate.Columns.Item(1).ColumnWidth = 30.666666666666668
$wsCreate.Columns.Item(2).ColumnWidth = 23.5

[void]$wsCreate.Activate()
[void]$wsCreate.Range("B2").Select()

# LoginTest ends up the active tab with A1:C2 selected (active cell C2).
[void]$wsLogin.Activate()
[void]$wsLogin.Range("A1:C2").Select()
